$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set header row (row 1)
$ws.Range("B1").Value = "N"
$ws.Range("C1").Value = "R2.scale_trim"
$ws.Range("D1").Value = "RMSE.scale_trim"
$ws.Range("E1").Value = "SMD.scale_trim"
$ws.Range("F1").Value = "SMD.scale_trim_round"
$ws.Range("G1").Value = "adj_agr.scale_trim_round"
$ws.Range("H1").Value = "corr.scale_trim"
$ws.Range("I1").Value = "exact_agr.scale_trim_round"
$ws.Range("J1").Value = "h_mean"
$ws.Range("K1").Value = "h_sd"
$ws.Range("L1").Value = "kappa.scale_trim_round"
$ws.Range("M1").Value = "sys_mean.scale_trim"
$ws.Range("N1").Value = "sys_mean.scale_trim_round"
$ws.Range("O1").Value = "sys_sd.scale_trim"
$ws.Range("P1").Value = "sys_sd.scale_trim_round"
$ws.Range("Q1").Value = "wtkappa.scale_trim_round"

# Set data rows
# Row 2
$ws.Range("A2").Value = "All data"
$ws.Range("B2").Value = 200
$ws.Range("C2").Value = 0.5493203316759327
$ws.Range("D2").Value = 0.6189327249996216
$ws.Range("E2").Value = 0.02324751973535285
$ws.Range("F2").Value = -0.01054416404922412
$ws.Range("G2").Value = 99
$ws.Range("H2").Value = 0.7801773732608226
$ws.Range("I2").Value = 64
$ws.Range("J2").Value = 3.5
$ws.Range("K2").Value = 0.9242680113386591
$ws.Range("L2").Value = 0.4694178334561533
$ws.Range("M2").Value = 3.521731815875496
$ws.Range("N2").Value = 3.49
$ws.Range("O2").Value = 0.9452175731243878
$ws.Range("P2").Value = 0.9719172370000901
$ws.Range("Q2").Value = 0.782122905027933

# Row 3
$ws.Range("A3").Value = "QUESTION_1"
$ws.Range("B3").Value = 40
$ws.Range("C3").Value = 0.5603717064330405
$ws.Range("D3").Value = 0.6112970223483144
$ws.Range("E3").Value = 0.01162334957758026
$ws.Range("F3").Value = -0.02636041012306077
$ws.Range("G3").Value = 97.5
$ws.Range("H3").Value = 0.7844718465545462
$ws.Range("I3").Value = 65
$ws.Range("J3").Value = 3.5
$ws.Range("K3").Value = 0.9336995618478525
$ws.Range("L3").Value = 0.4805194805194805
$ws.Range("M3").Value = 3.510865524397958
$ws.Range("N3").Value = 3.475
$ws.Range("O3").Value = 0.9512048837755814
$ws.Range("P3").Value = 1.012422836565829
$ws.Range("Q3").Value = 0.7702702702702703

# Row 4
$ws.Range("A4").Value = "QUESTION_2"
$ws.Range("B4").Value = 40
$ws.Range("C4").Value = 0.510955844326675
$ws.Range("D4").Value = 0.6447383440763595
$ws.Range("E4").Value = 0.0300768975131298
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 100
$ws.Range("H4").Value = 0.7610597455754698
$ws.Range("I4").Value = 60
$ws.Range("J4").Value = 3.5
$ws.Range("K4").Value = 0.9336995618478525
$ws.Range("L4").Value = 0.4165907019143118
$ws.Range("M4").Value = 3.52811592833568
$ws.Range("N4").Value = 3.5
$ws.Range("O4").Value = 0.9528707037260122
$ws.Range("P4").Value = 0.9870962335856491
$ws.Range("Q4").Value = 0.7777777777777778

# Row 5
$ws.Range("A5").Value = "QUESTION_3"
$ws.Range("B5").Value = 40
$ws.Range("C5").Value = 0.5461519181734076
$ws.Range("D5").Value = 0.6211045560552616
$ws.Range("E5").Value = 0.02836271304262913
$ws.Range("F5").Value = -0.05272082024612153
$ws.Range("G5").Value = 100
$ws.Range("H5").Value = 0.7759113437727831
$ws.Range("I5").Value = 70
$ws.Range("J5").Value = 3.5
$ws.Range("K5").Value = 0.9336995618478525
$ws.Range("L5").Value = 0.5596330275229358
$ws.Range("M5").Value = 3.526513506154147
$ws.Range("N5").Value = 3.45
$ws.Range("O5").Value = 0.9435584036333717
$ws.Range("P5").Value = 0.95943359359198
$ws.Range("Q5").Value = 0.8285714285714285

# Row 6
$ws.Range("A6").Value = "QUESTION_4"
$ws.Range("B6").Value = 40
$ws.Range("C6").Value = 0.5955885216100645
$ws.Range("D6").Value = 0.5863017624324911
$ws.Range("E6").Value = 0.02212661061487138
$ws.Range("F6").Value = -0.05272082024612153
$ws.Range("G6").Value = 100
$ws.Range("H6").Value = 0.8050680429782734
$ws.Range("I6").Value = 65
$ws.Range("J6").Value = 3.5
$ws.Range("K6").Value = 0.9336995618478525
$ws.Range("L6").Value = 0.4843462246777165
$ws.Range("M6").Value = 3.520683988369733
$ws.Range("N6").Value = 3.45
$ws.Range("O6").Value = 0.9646953505756862
$ws.Range("P6").Value = 1.01147265068163
$ws.Range("Q6").Value = 0.8108108108108109

# Row 7
$ws.Range("A7").Value = "QUESTION_5"
$ws.Range("B7").Value = 40
$ws.Range("C7").Value = 0.5335336678364757
$ws.Range("D7").Value = 0.62967958704328
$ws.Range("E7").Value = 0.02404802792855174
$ws.Range("F7").Value = 0.07908123036918277
$ws.Range("G7").Value = 97.5
$ws.Range("H7").Value = 0.7742921256011517
$ws.Range("I7").Value = 60
$ws.Range("J7").Value = 3.5
$ws.Range("K7").Value = 0.9336995618478525
$ws.Range("L7").Value = 0.4057567316620241
$ws.Range("M7").Value = 3.522480132119958
$ws.Range("N7").Value = 3.575
$ws.Range("O7").Value = 0.9617353297337229
$ws.Range("P7").Value = 0.9306049865682223
$ws.Range("Q7").Value = 0.7205882352941176
